$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.345980666666667
$ws.Range("H2").Value = 4.037942
$ws.Range("I2").Value = 0.4408299556445331
$ws.Range("J2").Value = 0.4408299556445331
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 11.94076733333333
$ws.Range("N2").Value = 35.82230199999999
$ws.Range("O2").Value = 0.351152845403141
$ws.Range("P2").Value = 0.351152845403141
$ws.Range("Q2").Value = 16.07204197583155
$ws.Range("R2").Value = 144.648377782484
$ws.Range("S2").Value = 0.1547986932635182
$ws.Range("T2").Value = 0.1547986932635182

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.345980666666667
$ws.Range("H3").Value = 4.037942
$ws.Range("I3").Value = 0.4408299556445331
$ws.Range("J3").Value = 0.4408299556445331
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.91523933333333
$ws.Range("N3").Value = 38.745718
$ws.Range("O3").Value = 0.3798100167568152
$ws.Range("P3").Value = 0.3798100167568153
$ws.Range("Q3").Value = 17.38366244803955
$ws.Range("R3").Value = 156.452962032356
$ws.Range("S3").Value = 0.1674316328402562
$ws.Range("T3").Value = 0.1674316328402563

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.345980666666667
$ws.Range("H4").Value = 4.037942
$ws.Range("I4").Value = 0.4408299556445331
$ws.Range("J4").Value = 0.4408299556445331
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08745933333333333
$ws.Range("N4").Value = 0.262378
$ws.Range("O4").Value = 0.002571994989913974
$ws.Range("P4").Value = 0.002571994989913974
$ws.Range("Q4").Value = 0.1177185717862222
$ws.Range("R4").Value = 1.059467146076
$ws.Range("S4").Value = 0.001133812437321738
$ws.Range("T4").Value = 0.001133812437321739

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.345980666666667
$ws.Range("H5").Value = 4.037942
$ws.Range("I5").Value = 0.4408299556445331
$ws.Range("J5").Value = 0.4408299556445331
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.061006666666666
$ws.Range("N5").Value = 27.18302
$ws.Range("O5").Value = 0.2664651428501298
$ws.Range("P5").Value = 0.2664651428501298
$ws.Range("Q5").Value = 12.19593979387111
$ws.Range("R5").Value = 109.76345814484
$ws.Range("S5").Value = 0.1174658171034369
$ws.Range("T5").Value = 0.1174658171034369

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Nppc"
$ws.Range("C6").Value = "Npr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.707307
$ws.Range("H6").Value = 5.121921
$ws.Range("I6").Value = 0.5591700443554668
$ws.Range("J6").Value = 0.5591700443554669
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 11.94076733333333
$ws.Range("N6").Value = 35.82230199999999
$ws.Range("O6").Value = 0.351152845403141
$ws.Range("P6").Value = 0.351152845403141
$ws.Range("Q6").Value = 20.38655565357133
$ws.Range("R6").Value = 183.4790008821419
$ws.Range("S6").Value = 0.1963541521396227
$ws.Range("T6").Value = 0.1963541521396227

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Nppc"
$ws.Range("C7").Value = "Npr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.707307
$ws.Range("H7").Value = 5.121921
$ws.Range("I7").Value = 0.5591700443554668
$ws.Range("J7").Value = 0.5591700443554669
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.91523933333333
$ws.Range("N7").Value = 38.745718
$ws.Range("O7").Value = 0.3798100167568152
$ws.Range("P7").Value = 0.3798100167568153
$ws.Range("Q7").Value = 22.05027852047533
$ws.Range("R7").Value = 198.4525066842779
$ws.Range("S7").Value = 0.212378383916559
$ws.Range("T7").Value = 0.2123783839165591

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nppc"
$ws.Range("C8").Value = "Npr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.707307
$ws.Range("H8").Value = 5.121921
$ws.Range("I8").Value = 0.5591700443554668
$ws.Range("J8").Value = 0.5591700443554669
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08745933333333333
$ws.Range("N8").Value = 0.262378
$ws.Range("O8").Value = 0.002571994989913974
$ws.Range("P8").Value = 0.002571994989913974
$ws.Range("Q8").Value = 0.1493199320153333
$ws.Range("R8").Value = 1.343879388138
$ws.Range("S8").Value = 0.001438182552592235
$ws.Range("T8").Value = 0.001438182552592236

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nppc"
$ws.Range("C9").Value = "Npr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.707307
$ws.Range("H9").Value = 5.121921
$ws.Range("I9").Value = 0.5591700443554668
$ws.Range("J9").Value = 0.5591700443554669
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.061006666666666
$ws.Range("N9").Value = 27.18302
$ws.Range("O9").Value = 0.2664651428501298
$ws.Range("P9").Value = 0.2664651428501298
$ws.Range("Q9").Value = 15.46992010904667
$ws.Range("R9").Value = 139.22928098142
$ws.Range("S9").Value = 0.1489993257466929
$ws.Range("T9").Value = 0.1489993257466929
